$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old Z1/AA1 header cells
$ws.Range("Z1").Clear()
$ws.Range("AA1").Clear()

# Row 2 - AccountCode
$ws.Range("Z2").Value = "Provides the Accounting Codes used by Finance & Accounts Department"
$ws.Range("AA2").Value = "Accounting movements having financial implications on the contract"

# Row 3 - AccountAmt
$ws.Range("Z3").Value = "Amount as per the transactions"

# Row 4 - SeqNo
$ws.Range("Z4").Value = "Seqence Number "

# Row 5 - GlSign
$ws.Range("Z5").Value = "Group Ledger Sign (Debit or Credit)"

$helpRange = $ws.Range("Z2:Z5")
$helpRange.Borders.LineStyle = 1
$helpRange.Borders.Weight = 2
$helpRange.Font.Size = 7.5
$helpRange.VerticalAlignment = -4108
$helpRange.WrapText = $true

$aa2Range = $ws.Range("AA2")
$aa2Range.Borders.LineStyle = 1
$aa2Range.Borders.Weight = 2
$aa2Range.Font.Size = 7.5
$aa2Range.VerticalAlignment = -4108
$aa2Range.WrapText = $true

Write-Host "done"
